$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:T10").ClearContents()
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Edn1"
$ws.Cells.Item(2,3).Value = "Ednrb"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 9.957023333333334
$ws.Cells.Item(2,8).Value = 29.87107
$ws.Cells.Item(2,9).Value = 0.7839926662698464
$ws.Cells.Item(2,10).Value = 0.7839926662698464
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 21.22137533333333
$ws.Cells.Item(2,14).Value = 63.664126
$ws.Cells.Item(2,15).Value = 0.2719819326156953
$ws.Cells.Item(2,16).Value = 0.2719819326156953
$ws.Cells.Item(2,17).Value = 211.3017293594245
$ws.Cells.Item(2,18).Value = 1901.71556423482
$ws.Cells.Item(2,19).Value = 0.2132318405286047
$ws.Cells.Item(2,20).Value = 0.2132318405286047
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Edn1"
$ws.Cells.Item(3,3).Value = "Ednrb"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 9.957023333333334
$ws.Cells.Item(3,8).Value = 29.87107
$ws.Cells.Item(3,9).Value = 0.7839926662698464
$ws.Cells.Item(3,10).Value = 0.7839926662698464
$ws.Cells.Item(3,11).Value = 1
$ws.Cells.Item(3,12).Value = 0.3333333333333333
$ws.Cells.Item(3,13).Value = 0.01207966666666667
$ws.Cells.Item(3,14).Value = 0.036239
$ws.Cells.Item(3,15).Value = 0.0001548180093772148
$ws.Cells.Item(3,16).Value = 0.0001548180093772148
$ws.Cells.Item(3,17).Value = 0.1202775228588889
$ws.Cells.Item(3,18).Value = 1.08249770573
$ws.Cells.Item(3,19).Value = 0.0001213761839582327
$ws.Cells.Item(3,20).Value = 0.0001213761839582327
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Edn1"
$ws.Cells.Item(4,3).Value = "Ednrb"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 9.957023333333334
$ws.Cells.Item(4,8).Value = 29.87107
$ws.Cells.Item(4,9).Value = 0.7839926662698464
$ws.Cells.Item(4,10).Value = 0.7839926662698464
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 56.79149
$ws.Cells.Item(4,14).Value = 170.37447
$ws.Cells.Item(4,15).Value = 0.7278632493749275
$ws.Cells.Item(4,16).Value = 0.7278632493749274
$ws.Cells.Item(4,17).Value = 565.4741910647667
$ws.Cells.Item(4,18).Value = 5089.2677195829
$ws.Cells.Item(4,19).Value = 0.5706394495572835
$ws.Cells.Item(4,20).Value = 0.5706394495572835
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Edn1"
$ws.Cells.Item(5,3).Value = "Ednrb"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.539481333333333
$ws.Cells.Item(5,8).Value = 7.618444
$ws.Cells.Item(5,9).Value = 0.1999528046497
$ws.Cells.Item(5,10).Value = 0.1999528046497
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 21.22137533333333
$ws.Cells.Item(5,14).Value = 63.664126
$ws.Cells.Item(5,15).Value = 0.2719819326156953
$ws.Cells.Item(5,16).Value = 0.2719819326156953
$ws.Cells.Item(5,17).Value = 53.89128652666044
$ws.Cells.Item(5,18).Value = 485.0215787399441
$ws.Cells.Item(5,19).Value = 0.05438355024055398
$ws.Cells.Item(5,20).Value = 0.054383550240554
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Edn1"
$ws.Cells.Item(6,3).Value = "Ednrb"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.539481333333333
$ws.Cells.Item(6,8).Value = 7.618444
$ws.Cells.Item(6,9).Value = 0.1999528046497
$ws.Cells.Item(6,10).Value = 0.1999528046497
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.01207966666666667
$ws.Cells.Item(6,14).Value = 0.036239
$ws.Cells.Item(6,15).Value = 0.0001548180093772148
$ws.Cells.Item(6,16).Value = 0.0001548180093772148
$ws.Cells.Item(6,17).Value = 0.03067608801288889
$ws.Cells.Item(6,18).Value = 0.276084792116
$ws.Cells.Item(6,19).Value = [double]"3.095629518525765E-05"
$ws.Cells.Item(6,20).Value = [double]"3.095629518525765E-05"
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Edn1"
$ws.Cells.Item(7,3).Value = "Ednrb"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.539481333333333
$ws.Cells.Item(7,8).Value = 7.618444
$ws.Cells.Item(7,9).Value = 0.1999528046497
$ws.Cells.Item(7,10).Value = 0.1999528046497
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 56.79149
$ws.Cells.Item(7,14).Value = 170.37447
$ws.Cells.Item(7,15).Value = 0.7278632493749275
$ws.Cells.Item(7,16).Value = 0.7278632493749274
$ws.Cells.Item(7,17).Value = 144.2209287471867
$ws.Cells.Item(7,18).Value = 1297.98835872468
$ws.Cells.Item(7,19).Value = 0.1455382981139607
$ws.Cells.Item(7,20).Value = 0.1455382981139607
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Edn1"
$ws.Cells.Item(8,3).Value = "Ednrb"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.203899
$ws.Cells.Item(8,8).Value = 0.611697
$ws.Cells.Item(8,9).Value = 0.01605452908045364
$ws.Cells.Item(8,10).Value = 0.01605452908045364
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 21.22137533333333
$ws.Cells.Item(8,14).Value = 63.664126
$ws.Cells.Item(8,15).Value = 0.2719819326156953
$ws.Cells.Item(8,16).Value = 0.2719819326156953
$ws.Cells.Item(8,17).Value = 4.327017209091334
$ws.Cells.Item(8,18).Value = 38.943154881822
$ws.Cells.Item(8,19).Value = 0.004366541846536663
$ws.Cells.Item(8,20).Value = 0.004366541846536663
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Edn1"
$ws.Cells.Item(9,3).Value = "Ednrb"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.203899
$ws.Cells.Item(9,8).Value = 0.611697
$ws.Cells.Item(9,9).Value = 0.01605452908045364
$ws.Cells.Item(9,10).Value = 0.01605452908045364
$ws.Cells.Item(9,11).Value = 1
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.01207966666666667
$ws.Cells.Item(9,14).Value = 0.036239
$ws.Cells.Item(9,15).Value = 0.0001548180093772148
$ws.Cells.Item(9,16).Value = 0.0001548180093772148
$ws.Cells.Item(9,17).Value = 0.002463031953666667
$ws.Cells.Item(9,18).Value = 0.022167287583
$ws.Cells.Item(9,19).Value = [double]"2.485530233724439E-06"
$ws.Cells.Item(9,20).Value = [double]"2.485530233724439E-06"
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Edn1"
$ws.Cells.Item(10,3).Value = "Ednrb"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.203899
$ws.Cells.Item(10,8).Value = 0.611697
$ws.Cells.Item(10,9).Value = 0.01605452908045364
$ws.Cells.Item(10,10).Value = 0.01605452908045364
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 56.79149
$ws.Cells.Item(10,14).Value = 170.37447
$ws.Cells.Item(10,15).Value = 0.7278632493749275
$ws.Cells.Item(10,16).Value = 0.7278632493749274
$ws.Cells.Item(10,17).Value = 11.57972801951
$ws.Cells.Item(10,18).Value = 104.21755217559
$ws.Cells.Item(10,19).Value = 0.01168550170368325
$ws.Cells.Item(10,20).Value = 0.01168550170368325
Write-Output "done"
